# Generate Report for Handoff
# Update the localization-status workbook: the handback markdown file was
# renamed from a6eef75e-519e-4e4a-a75b-96fbd026565b.md to
# e7ac1c55-732b-428d-97fa-330a03d2eb37.md, a new xliff handoff round was
# generated (new content hash + new timestamps), and the report text was
# regenerated accordingly (hyperlink targets are left pointing at the
# original commit URL - only the cell text / hyperlink display text is
# refreshed).

$wb = $excel.ActiveWorkbook

$oldGuid = "a6eef75e-519e-4e4a-a75b-96fbd026565b"
$newGuid = "e7ac1c55-732b-428d-97fa-330a03d2eb37"

$oldHash = "60fb928e4e669cbe583bc2f449a61f1662fbe85d"
$newHash = "79e5fc2712d81f9d4ba3843784a8f99c661e8e32"

$handoffUrl = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/1a9ec81ab564503d8bf8c92c9f35d36114d45b5e/e2e/$oldGuid.md"

function Set-HyperlinkCell($ws, $cellAddr, $displayText) {
    $rng = $ws.Range($cellAddr)
    $url = $rng.Hyperlinks.Item(1).Address
    if ([string]::IsNullOrEmpty($url)) { $url = $handoffUrl }
    $rng.Hyperlinks.Delete()
    $rng.Value = $displayText
    $ws.Hyperlinks.Add($rng, $url, "", "", $displayText)
}

# ---- Sheet "Overview" ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
Set-HyperlinkCell $wsOverview "B2" "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2017-02-22 07:47:28"

# ---- Sheet "zh-cn" ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HyperlinkCell $wsZhCn "A2" "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2017-02-22 07:47:11"

# ---- Sheet "de-de" ----
$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HyperlinkCell $wsDeDe "A2" "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2017-02-22 07:47:28"
